$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Supermercado" rows appended below the existing data (rows 13-22).
# Columns: A categoria, B nome, C latitude, D longitude, E icone, G descricao
$rows = @(
    @{ Row=13; Nome="Tesco Express";      Lat=55.943367940442798;  Lon=-3.20957015359396;   Desc="Horário: 06:00–23:00" },
    @{ Row=14; Nome="Tesco Express";      Lat=55.9436566678961;    Lon=-3.2093410629952701; Desc="Horário: 06:00–23:00" },
    @{ Row=15; Nome="Tesco Express";      Lat=55.948082610613298;  Lon=-3.1860243139898499; Desc="Horário: 06:00–23:00" },
    @{ Row=16; Nome="Tesco Express";      Lat=55.950615650317602;  Lon=-3.1784637130764302; Desc="Horário: 06:00–23:00" },
    @{ Row=17; Nome="Tesco Express";      Lat=55.957352691313702;  Lon=-3.18693992108699;   Desc="Horário: 06:00–23:00" },
    @{ Row=18; Nome="Tesco Express";      Lat=55.946519798700301;  Lon=-3.2225433733526998; Desc="Horário: 06:00–23:00" },
    @{ Row=19; Nome="Poundland";          Lat=55.945138507949402;  Lon=3.2056257210258399;  Desc="Horário: 08:00–20:00" },
    @{ Row=20; Nome="Sainsbury's Local";  Lat=55.9463749104114;    Lon=-3.2012609818517901; Desc="Horário: 07:00–23:00" },
    @{ Row=21; Nome="Sainsbury's Local";  Lat=55.949017986035301;  Lon=-3.1870979385102198; Desc="Horário: 07:00–23:00" },
    @{ Row=22; Nome="LIDL";               Lat=55.9458654868656;    Lon=-3.1844826429666102; Desc="Horário: 07:00–22:00" }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("A$i").Value = "Supermercado"
    $ws.Range("A$i").HorizontalAlignment = -4131
    $ws.Range("E$i").Value = "supermarket.png"
    $ws.Range("B$i").Value = $r.Nome
    $ws.Range("C$i").Value = $r.Lat
    $ws.Range("D$i").Value = $r.Lon
    $ws.Range("G$i").Value = $r.Desc
}

# Match the selection left behind after entering the new rows.
$ws.Range("A13:XFD13").Select()
